# fix of A/S problem
# Insert a new "Take away" slide right before the final "Q&A" slide by
# duplicating the preceding "Demo" slide (slide 26) and updating its text.

$p = $ppt.ActivePresentation

# Slide 26 ("Demo") is duplicated; PowerPoint places the duplicate
# immediately after it, i.e. right before the final "Q&A" slide (27).
$demoSlide = $p.Slides.Item(26)
$newSlide = $demoSlide.Duplicate()
$newSlide = $p.Slides.Item(27)

# Update the title placeholder text.
$newSlide.Shapes.Item("Title 1").TextFrame.TextRange.Text = "Take away"

# Update the body placeholder text.
$newSlide.Shapes.Item("Text Placeholder 2").TextFrame.TextRange.Text = "Model your data " + [char]8220 + "more" + [char]8221 + " like what is represents."
